# Applies the "Updated cryptos list" data refresh described by the diff.
# The workbook stores prices/volumes as plain text in columns D and E.
# Values that look like plain decimal numbers (single dot) are written with a
# leading apostrophe so Excel keeps them as text instead of auto-converting
# them to floating point numbers (this matches the original text content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "37.546.21"
$ws.Range("D3").Value = "2.059.97"
$ws.Range("D5").Value = "'253.73"
$ws.Range("D6").Value = "'0.654"
$ws.Range("D7").Value = "'67.54"
$ws.Range("D9").Value = "'0.396"
$ws.Range("D10").Value = "'59.93"
$ws.Range("D11").Value = "'0.0772"
$ws.Range("D13").Value = "'0.941"
$ws.Range("D14").Value = "'24.08"
$ws.Range("D15").Value = "'15.06"
$ws.Range("D16").Value = "2.365.55"
$ws.Range("D17").Value = "'5.71"
$ws.Range("D18").Value = "2.059.04"
$ws.Range("D19").Value = "37.481.94"
$ws.Range("D20").Value = "'73.82"
$ws.Range("D21").Value = "0.0₃0881"
$ws.Range("D22").Value = "'5.51"
$ws.Range("D23").Value = "'240.97"
$ws.Range("D24").Value = "'2.71"
$ws.Range("D26").Value = "'2.46"
$ws.Range("D27").Value = "'10.09"
$ws.Range("D33").Value = "'1.23"
$ws.Range("D34").Value = "'4.75"
$ws.Range("D37").Value = "'6.27"
$ws.Range("D40").Value = "'3.13"
$ws.Range("D42").Value = "'18.45"
$ws.Range("D47").Value = "'97.75"
$ws.Range("D48").Value = "'8.03"
$ws.Range("D49").Value = "1.421.30"
$ws.Range("D51").Value = "'3.80"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +5.32%  "
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +3.09%  "
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("E7").Value = "  +15.61%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +9.70%  "
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +4.76%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("E14").Value = "  +29.56%  "
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("E17").Value = "  +7.97%  "
$ws.Range("E18").Value = "  +3.83%  "
$ws.Range("E19").Value = "  +5.38%  "
$ws.Range("E20").Value = "  +3.21%  "
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("E22").Value = "  +5.37%  "
$ws.Range("E23").Value = "  +3.40%  "
$ws.Range("E24").Value = "  +4.84%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  +8.09%  "
$ws.Range("E27").Value = "  +10.11%  "
$ws.Range("E30").Value = "  +4.99%  "
$ws.Range("E31").Value = "  +3.30%  "
$ws.Range("E32").Value = "  +7.43%  "
$ws.Range("E33").Value = "  +10.47%  "
$ws.Range("E34").Value = "  +8.62%  "
$ws.Range("E35").Value = "  +6.41%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  +15.08%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("E40").Value = "  +37.57%  "
$ws.Range("E41").Value = "  +14.27%  "
$ws.Range("E42").Value = "  +14.10%  "
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("E44").Value = "  +6.70%  "
$ws.Range("E45").Value = "  +3.65%  "
$ws.Range("E46").Value = "  +5.71%  "
$ws.Range("E47").Value = "  +4.59%  "
$ws.Range("E48").Value = "  +3.21%  "
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("E51").Value = "  +10.11%  "

# --- Rows 28/29: Monero and Kaspa swapped ranking positions ---
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.140"
$ws.Range("E28").Value = "  +46.11%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'162.83"
$ws.Range("E29").Value = "  -1.33%  "

Write-Host "cryptos list updated"
